# Generate Report for Handback
# Rename the two tracked e2e markdown files (and their derived xliff hashes),
# and refresh the timestamps recorded for the zh-cn / de-de handback rows.
#
# Old file 1: 1aba61aa-d259-415f-9447-4f8db1ef9e15.md
# New file 1: e99a8093-d901-4ade-b759-d1188cbcd08a.md
#
# Old file 2: 1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md
# New file 2: ffff81fb709f-8978-4d75-8a9a-6197e3cad16a.md
#
# Old xliff hash: 548d2abe1ab53c22add390a263f740f8913a8fea / 0537f948374ccd930d7dfab2b0d917ab8642c0b5
# New xliff hash (shared by both rows now): 354ab020ab03dc3e85a322a98dfd26d154bc1c11

$wb = $excel.ActiveWorkbook

$oldName1 = "1aba61aa-d259-415f-9447-4f8db1ef9e15.md"
$newName1 = "e99a8093-d901-4ade-b759-d1188cbcd08a.md"
$oldName2 = "1fe2c9e0-7f52-4499-a4c6-0e4c1466dce6.md"
$newName2 = "ffff81fb709f-8978-4d75-8a9a-6197e3cad16a.md"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newName1
$ws1.Range("B2").Value = "e2e\" + $newName1
$ws1.Range("G2").Value = "2016-08-29 03:02:30"

$ws1.Range("A3").Value = $newName2
$ws1.Range("B3").Value = "e2e\" + $newName2
$ws1.Range("G3").Value = "2016-08-29 03:02:30"

# Hyperlinks in this engine can only be refreshed by deleting every
# hyperlink on the sheet and re-adding them - per-item Address/TextToDisplay
# edits silently create a duplicate, orphaned <hyperlink> element instead of
# updating in place. The link *targets* are untouched upstream, so re-use
# the original Address values and only change the displayed text.
$rngB2 = $ws1.Range("B2")
$rngB3 = $ws1.Range("B3")
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($rngB2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName1, [Type]::Missing, [Type]::Missing, "e2e\" + $newName1)
$ws1.Hyperlinks.Add($rngB3, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName2, [Type]::Missing, [Type]::Missing, "e2e\" + $newName2)

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newName1
$ws2.Range("I2").Value = $newName1
$ws2.Range("A3").Value = $newName2
$ws2.Range("I3").Value = $newName2

# Correspond Handoff File / Correspond Handback File / their datetimes:
# both rows now point at the same regenerated xliff.
$newXlf1 = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.zh-cn.xlf"

$ws2.Range("G2").Value = $newXlf1
$ws2.Range("H2").Value = "2016-08-29 03:02:24"
$ws2.Range("J2").Value = $newXlf1
$ws2.Range("K2").Value = "2016-08-29 03:02:41"

$ws2.Range("G3").Value = $newXlf1
$ws2.Range("H3").Value = "2016-08-29 03:02:24"
$ws2.Range("J3").Value = $newXlf1
$ws2.Range("K3").Value = "2016-08-29 03:02:41"

$rngA2b = $ws2.Range("A2")
$rngI2b = $ws2.Range("I2")
$rngA3b = $ws2.Range("A3")
$rngI3b = $ws2.Range("I3")
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($rngA2b, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName1, [Type]::Missing, [Type]::Missing, $newName1)
$ws2.Hyperlinks.Add($rngI2b, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/874e02670c277fc821649de3ae4a43c0ea67e62b/e2e/" + $oldName1, [Type]::Missing, [Type]::Missing, $newName1)
$ws2.Hyperlinks.Add($rngA3b, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName2, [Type]::Missing, [Type]::Missing, $newName2)
$ws2.Hyperlinks.Add($rngI3b, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/874e02670c277fc821649de3ae4a43c0ea67e62b/e2e/" + $oldName2, [Type]::Missing, [Type]::Missing, $newName2)

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newName1
$ws3.Range("I2").Value = $newName1
$ws3.Range("A3").Value = $newName2
$ws3.Range("I3").Value = $newName2

$newXlf2 = "e99a8093-d901-4ade-b759-d1188cbcd08a.354ab020ab03dc3e85a322a98dfd26d154bc1c11.de-de.xlf"

$ws3.Range("G2").Value = $newXlf2
$ws3.Range("H2").Value = "2016-08-29 03:02:30"
$ws3.Range("J2").Value = $newXlf2
$ws3.Range("K2").Value = "2016-08-29 03:02:48"

$ws3.Range("G3").Value = $newXlf2
$ws3.Range("H3").Value = "2016-08-29 03:02:30"
$ws3.Range("J3").Value = $newXlf2
$ws3.Range("K3").Value = "2016-08-29 03:02:48"

$rngA2c = $ws3.Range("A2")
$rngI2c = $ws3.Range("I2")
$rngA3c = $ws3.Range("A3")
$rngI3c = $ws3.Range("I3")
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($rngA2c, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName1, [Type]::Missing, [Type]::Missing, $newName1)
$ws3.Hyperlinks.Add($rngI2c, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/78723b2aea5b0d144ab946809f5befc72dd9629f/e2e/" + $oldName1, [Type]::Missing, [Type]::Missing, $newName1)
$ws3.Hyperlinks.Add($rngA3c, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0cacdb166865b23e85868a85cb7be2548a36fb23/e2e/" + $oldName2, [Type]::Missing, [Type]::Missing, $newName2)
$ws3.Hyperlinks.Add($rngI3c, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/78723b2aea5b0d144ab946809f5befc72dd9629f/e2e/" + $oldName2, [Type]::Missing, [Type]::Missing, $newName2)
